# Auto-generated edit script: update market price columns (H-N) across multiple Tonberry Profit sheets
$wb = $excel.ActiveWorkbook

# ALC!row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6123.625
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4516

# ALC!row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1396.3334
$ws.Range("J125").Value = 1333.375
$ws.Range("L125").Value = 12000.375
$ws.Range("N125").Value = -16920.375

# ARM!row 23
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 30951.25
$ws.Range("J23").Value = 17933
$ws.Range("L23").Value = 17933
$ws.Range("N23").Value = -18451

# ARM!row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4124.5093
$ws.Range("I32").Value = 2705.4358
$ws.Range("J32").Value = 7583.5
$ws.Range("K32").Value = 2705.4358
$ws.Range("L32").Value = 7583.5
$ws.Range("M32").Value = -2418.4358
$ws.Range("N32").Value = -8157.5

# ARM!row 34
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 73028
$ws.Range("J34").Value = 73028
$ws.Range("L34").Value = 73028
$ws.Range("N34").Value = -73570

# ARM!row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1605.5
$ws.Range("I45").Value = 1174.75
$ws.Range("J45").Value = 1892.6666
$ws.Range("K45").Value = 1174.75
$ws.Range("L45").Value = 1892.6666
$ws.Range("M45").Value = -797.75
$ws.Range("N45").Value = -2646.6666

# ARM!row 47
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

# ARM!row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2749.2083
$ws.Range("I61").Value = 2244.8
$ws.Range("J61").Value = 5271.25
$ws.Range("K61").Value = 2244.8
$ws.Range("L61").Value = 5271.25
$ws.Range("M61").Value = -2032.8
$ws.Range("N61").Value = -5695.25

# ARM!row 109
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 99997.5
$ws.Range("J109").Value = 99997.5
$ws.Range("L109").Value = 99997.5
$ws.Range("N109").Value = -102771.5

# ARM!row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1575.4445
$ws.Range("I122").Value = 1293.8
$ws.Range("K122").Value = 3881.4
$ws.Range("M122").Value = -1431.4

# ARM!row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2749.2083
$ws.Range("I136").Value = 2244.8
$ws.Range("J136").Value = 5271.25
$ws.Range("K136").Value = 6734.400000000001
$ws.Range("L136").Value = 15813.75
$ws.Range("M136").Value = -4184.400000000001
$ws.Range("N136").Value = -20913.75

# BSM!row 9
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

# BSM!row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1901.2
$ws.Range("I107").Value = 2228.5715
$ws.Range("K107").Value = 2228.5715
$ws.Range("M107").Value = -308.5715

# CRP!row 13
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 1498
$ws.Range("J13").Value = 1498
$ws.Range("L13").Value = 1498
$ws.Range("N13").Value = -1776

# CRP!row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 618.53845
$ws.Range("I22").Value = 249.16667
$ws.Range("J22").Value = 935.1429000000001
$ws.Range("K22").Value = 249.16667
$ws.Range("L22").Value = 935.1429000000001
$ws.Range("M22").Value = 100.83333
$ws.Range("N22").Value = -1635.1429

# CRP!row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()

# CRP!row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1891570.4
$ws.Range("I58").Value = 2718383
$ws.Range("J58").Value = 1712.8572
$ws.Range("K58").Value = 2718383
$ws.Range("L58").Value = 1712.8572
$ws.Range("M58").Value = -2718180
$ws.Range("N58").Value = -2118.8572

# CRP!row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1358
$ws.Range("I122").Value = 1301.9445
$ws.Range("J122").Value = 1526.1666
$ws.Range("K122").Value = 3905.8335
$ws.Range("L122").Value = 4578.4998
$ws.Range("M122").Value = -1455.8335
$ws.Range("N122").Value = -9478.4998

# CRP!row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1891570.4
$ws.Range("I136").Value = 2718383
$ws.Range("J136").Value = 1712.8572
$ws.Range("K136").Value = 8155149
$ws.Range("L136").Value = 5138.571599999999
$ws.Range("M136").Value = -8152599
$ws.Range("N136").Value = -10238.5716

# CUL!row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 27751
$ws.Range("I55").Value = 35334.668
$ws.Range("J55").Value = 5000
$ws.Range("K55").Value = 106004.004
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = -105827.004
$ws.Range("N55").Value = -15354

# CUL!row 110
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 3000
$ws.Range("I110").Value = 3000
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 9000
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -4910
$ws.Range("N110").ClearContents()

# GSM!row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2612.4285
$ws.Range("I122").Value = 2349.75
$ws.Range("J122").Value = 2962.6667
$ws.Range("K122").Value = 7049.25
$ws.Range("L122").Value = 8888.000100000001
$ws.Range("M122").Value = -4599.25
$ws.Range("N122").Value = -13788.0001

# LTW!row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3465.0908
$ws.Range("I7").Value = 1708.0625
$ws.Range("J7").Value = 8150.5
$ws.Range("K7").Value = 1708.0625
$ws.Range("L7").Value = 8150.5
$ws.Range("M7").Value = -1596.0625
$ws.Range("N7").Value = -8374.5

# LTW!row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2051.85
$ws.Range("J16").Value = 552.875
$ws.Range("L16").Value = 552.875
$ws.Range("N16").Value = -892.875

# LTW!row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2000
$ws.Range("J22").Value = 3000
$ws.Range("L22").Value = 3000
$ws.Range("N22").Value = -3590

# LTW!row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2000
$ws.Range("J27").Value = 3000
$ws.Range("L27").Value = 3000
$ws.Range("N27").Value = -3214

# LTW!row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2258.7
$ws.Range("J61").Value = 2460.6
$ws.Range("L61").Value = 2460.6
$ws.Range("N61").Value = -2864.6

# LTW!row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2258.7
$ws.Range("J113").Value = 2460.6
$ws.Range("L113").Value = 2460.6
$ws.Range("N113").Value = -6800.6

# LTW!row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3720.5
$ws.Range("I122").Value = 3399.111
$ws.Range("K122").Value = 10197.333
$ws.Range("M122").Value = -7747.332999999999

# LTW!row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3465.0908
$ws.Range("I126").Value = 1708.0625
$ws.Range("J126").Value = 8150.5
$ws.Range("K126").Value = 5124.1875
$ws.Range("L126").Value = 24451.5
$ws.Range("M126").Value = -2654.1875
$ws.Range("N126").Value = -29391.5

